$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.431.97'
$ws.Range('E2').Value = '  +3.67%  '
$ws.Range('D3').Value = '2.059.35'
$ws.Range('E3').Value = '  +5.36%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('E6').Value = '  +3.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.84'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.84%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +4.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.16'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.08%  '
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('E12').Value = '  +4.16%  '
$ws.Range('D13').Value = '2.363.45'
$ws.Range('E13').Value = '  +5.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.24'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.84'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.773'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.18'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.56%  '
$ws.Range('D18').Value = '2.062.06'
$ws.Range('E18').Value = '  +5.70%  '
$ws.Range('D19').Value = '37.607.00'
$ws.Range('E19').Value = '  +4.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.18'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +24.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.71'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.44%  '
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '224.68'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.67%  '
$ws.Range('E26').Value = '  +3.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.71'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.88'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.41'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.75%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.26'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.30%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.126'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.119'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.50'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +5.62%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0629'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.74%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.60'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +15.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.46'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.18%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.79'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.33'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +6.65%  '
$ws.Range('E40').Value = '  +13.66%  '
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.74'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +37.17%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.98'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0965'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +11.09%  '
$ws.Range('D44').Value = '1.473.21'
$ws.Range('E44').Value = '  +5.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '95.47'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +11.19%  '
$ws.Range('E46').Value = '  +6.89%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.14'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +7.38%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.97'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +10.18%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.02'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.07%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.25'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +8.99%  '
$ws.Range('E51').Value = '  +4.05%  '
